# "multi browser implementation for chrome and firefox and accounting cash"
#
# Two semantic changes to the ProductLoanInput sheet:
#  1. The "Currency" row's label/value are re-cased / cleaned up:
#       A6: "Currency"      -> "currency"
#       B6: "US Dollar "    -> "US Dollar"   (drop the trailing space)
#  2. A new "accrualperiodic" (accounting: cash/accrual) flag row is
#     inserted right after "maximumallowedoutstandingbalancefortranchloan"
#     (i.e. before the "fundsource" accounting block), set to "checked".
#
# The stray, empty, formatted column-C cells left over in column C (rows
# 5, 6, 12) are cleaned up at the same time, which is why the sheet's used
# range shrinks from A1:C40 down to A1:B41.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# --- 1. Fix up the currency label / value on row 6 --------------------
$ws.Cells.Item(6, 1).Value = "currency"
$ws.Cells.Item(6, 2).Value = "US Dollar"

# --- 2. Drop the leftover empty column C ------------------------------
$ws.Columns.Item(3).Clear()

# --- 3. Insert the new "accrualperiodic" row before "fundsource" ------
$ws.Rows.Item(29).Insert()

$ws.Cells.Item(29, 1).Value = "accrualperiodic"
$ws.Cells.Item(29, 2).Value = "checked"

# Match formatting of the neighbouring tranche-section rows (grey label
# fill, green checked-value fill, wrapped Arial 10pt) for the new row.
$labelCell = $ws.Cells.Item(29, 1)
$labelCell.Interior.Color = $ws.Cells.Item(26, 1).Interior.Color
$labelCell.Font.Name = $ws.Cells.Item(26, 1).Font.Name
$labelCell.Font.Size = $ws.Cells.Item(26, 1).Font.Size
$labelCell.WrapText = $true

$valueCell = $ws.Cells.Item(29, 2)
$valueCell.Interior.Color = $ws.Cells.Item(26, 1).Interior.Color
$valueCell.Font.Name = $ws.Cells.Item(26, 1).Font.Name
$valueCell.Font.Size = $ws.Cells.Item(26, 1).Font.Size
$valueCell.WrapText = $true

# --- 4. Restore the selection to the edited currency row --------------
$ws.Range("A6:B6").Select()
